$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New "Progreso Día 5" block (rows 39-42), mirroring the layout of the
#     previous "Progreso Día 4" block (rows 34-37) ---

# Merge the title row first (on still-blank cells) so that the later
# format-copy doesn't get its borders split across the merged region.
$ws.Range("A39:E39").Merge() | Out-Null

# Copy formatting only (keeps style ids + avoids clobbering the row heights
# we set further below) from the equivalent rows of the prior section.
$ws.Range("A34:E34").Copy() | Out-Null
$ws.Range("A39:E39").PasteSpecial(-4122) | Out-Null

$ws.Range("A35:E35").Copy() | Out-Null
$ws.Range("A40:E40").PasteSpecial(-4122) | Out-Null

$ws.Range("A36:E36").Copy() | Out-Null
$ws.Range("A41:E41").PasteSpecial(-4122) | Out-Null

$ws.Range("A37:E37").Copy() | Out-Null
$ws.Range("A42:E42").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Section title (row 39, already merged across A:E above)
$ws.Range("A39").Value = "Progreso Día 5"

# Column headers (row 40)
$ws.Range("A40").Value = "Fecha"
$ws.Range("B40").Value = "Tarea"
$ws.Range("C40").Value = "Descripción del Progreso"
$ws.Range("D40").Value = "Archivos Modificados"
$ws.Range("E40").Value = "Observaciones"

# Row 41 - first entry (E41 is filled in further down, after row 42, to
# match the original authoring/shared-string order)
$ws.Range("A41").Value = 45496
$ws.Range("B41").Value = "Códigos Ordenados"
$ws.Range("C41").Value = "Se ordenaron los códigos con la extensión Prettier"
$ws.Range("D41").Value = "Todos los JS y JSX"

# Row 42 - second entry
$ws.Range("A42").Value = 45496
$ws.Range("B42").Value = "Eliminación de Línea de Código Innecesaria"
$ws.Range("C42").Value = "Se eliminó línea de código (const apiCountries = await fetchCountries();) que estaba sin uso debido a actualizaciones anteriores"
$ws.Range("D42").Value = "src/components/EditCountry.jsx"
$ws.Range("E42").Value = "Anten consola advertía de un riesgo, ya no lo hace después de este arreglo."

# Back to E41, filled last
$ws.Range("E41").Value = "Código más ordenado, bonito y fácil de entender."

# Row heights to match the rest of the table (auto-fit-like sizing)
$ws.Rows.Item(39).RowHeight = 15.75
$ws.Rows.Item(41).RowHeight = 60
$ws.Rows.Item(42).RowHeight = 75

# Scroll / select so the new rows are in view, like after manual entry
$excel.ActiveWindow.ScrollRow = 37
$excel.Goto($ws.Range("C42"), $true)
